$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 12
$ws_ALC.Range("H12").Value = 120
$ws_ALC.Range("I12").Value = 120
$ws_ALC.Range("J12").Value = 0
$ws_ALC.Range("K12").Value = 120
$ws_ALC.Range("L12").Value = 0
$ws_ALC.Range("M12").Value = 50
$ws_ALC.Range("N12").ClearContents()

# ALC row 113
$ws_ALC.Range("H113").Value = 14207.272
$ws_ALC.Range("I113").Value = 3793.3333
$ws_ALC.Range("J113").Value = 18112.5
$ws_ALC.Range("K113").Value = 3793.3333
$ws_ALC.Range("L113").Value = 18112.5
$ws_ALC.Range("M113").Value = -539.3332999999998
$ws_ALC.Range("N113").Value = -24620.5

# ALC row 114
$ws_ALC.Range("H114").Value = 39888
$ws_ALC.Range("J114").Value = 39888
$ws_ALC.Range("L114").Value = 39888
$ws_ALC.Range("N114").Value = -48566

# ALC row 121
$ws_ALC.Range("H121").Value = 3101.111
$ws_ALC.Range("I121").Value = 865
$ws_ALC.Range("J121").Value = 4219.1665
$ws_ALC.Range("K121").Value = 2595
$ws_ALC.Range("L121").Value = 12657.4995
$ws_ALC.Range("M121").Value = -848
$ws_ALC.Range("N121").Value = -16151.4995

# ALC row 135
$ws_ALC.Range("H135").Value = 883.05
$ws_ALC.Range("I135").Value = 913.4167
$ws_ALC.Range("J135").Value = 837.5
$ws_ALC.Range("K135").Value = 8220.7503
$ws_ALC.Range("L135").Value = 7537.5
$ws_ALC.Range("M135").Value = -5685.7503
$ws_ALC.Range("N135").Value = -12607.5

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 97
$ws_ARM.Range("H97").Value = 1335.6
$ws_ARM.Range("I97").Value = 1076.3914
$ws_ARM.Range("J97").Value = 2187.2856
$ws_ARM.Range("K97").Value = 1076.3914
$ws_ARM.Range("L97").Value = 2187.2856
$ws_ARM.Range("M97").Value = -580.3914
$ws_ARM.Range("N97").Value = -3179.2856

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws_CRP.Range("H7").Value = 804.7059
$ws_CRP.Range("I7").Value = 1518.1428
$ws_CRP.Range("J7").Value = 305.3
$ws_CRP.Range("K7").Value = 1518.1428
$ws_CRP.Range("L7").Value = 305.3
$ws_CRP.Range("M7").Value = -1405.1428
$ws_CRP.Range("N7").Value = -531.3

# CRP row 31
$ws_CRP.Range("H31").Value = 9263553
$ws_CRP.Range("I31").Value = 11112597
$ws_CRP.Range("J31").Value = 18333.334
$ws_CRP.Range("K31").Value = 11112597
$ws_CRP.Range("L31").Value = 18333.334
$ws_CRP.Range("M31").Value = -11112302
$ws_CRP.Range("N31").Value = -18923.334

# CRP row 34
$ws_CRP.Range("H34").Value = 9263553
$ws_CRP.Range("I34").Value = 11112597
$ws_CRP.Range("J34").Value = 18333.334
$ws_CRP.Range("K34").Value = 11112597
$ws_CRP.Range("L34").Value = 18333.334
$ws_CRP.Range("M34").Value = -11112395
$ws_CRP.Range("N34").Value = -18737.334

# CRP row 58
$ws_CRP.Range("H58").Value = 2375.652
$ws_CRP.Range("I58").Value = 1200.5454
$ws_CRP.Range("J58").Value = 3452.8333
$ws_CRP.Range("K58").Value = 1200.5454
$ws_CRP.Range("L58").Value = 3452.8333
$ws_CRP.Range("M58").Value = -997.5454
$ws_CRP.Range("N58").Value = -3858.8333

# CRP row 132
$ws_CRP.Range("H132").Value = 3075.6667
$ws_CRP.Range("I132").Value = 1447
$ws_CRP.Range("K132").Value = 4341
$ws_CRP.Range("M132").Value = -1811

# CRP row 136
$ws_CRP.Range("H136").Value = 2375.652
$ws_CRP.Range("I136").Value = 1200.5454
$ws_CRP.Range("J136").Value = 3452.8333
$ws_CRP.Range("K136").Value = 3601.6362
$ws_CRP.Range("L136").Value = 10358.4999
$ws_CRP.Range("M136").Value = -1051.6362
$ws_CRP.Range("N136").Value = -15458.4999

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 2
$ws_CUL.Range("H2").Value = 106.85714
$ws_CUL.Range("I2").Value = 40.333332
$ws_CUL.Range("J2").Value = 156.75
$ws_CUL.Range("K2").Value = 241.999992
$ws_CUL.Range("L2").Value = 940.5
$ws_CUL.Range("M2").Value = -128.999992
$ws_CUL.Range("N2").Value = -1166.5

# CUL row 7
$ws_CUL.Range("H7").Value = 146.66667
$ws_CUL.Range("I7").Value = 146.66667
$ws_CUL.Range("J7").Value = 0
$ws_CUL.Range("K7").Value = 440.00001
$ws_CUL.Range("L7").Value = 0
$ws_CUL.Range("M7").Value = -328.00001
$ws_CUL.Range("N7").ClearContents()

# CUL row 11
$ws_CUL.Range("H11").Value = 78.09090999999999
$ws_CUL.Range("I11").Value = 78.09090999999999
$ws_CUL.Range("K11").Value = 234.27273
$ws_CUL.Range("M11").Value = -94.27272999999997

# CUL row 19
$ws_CUL.Range("H19").Value = 2466.6667
$ws_CUL.Range("I19").Value = 0
$ws_CUL.Range("J19").Value = 2466.6667
$ws_CUL.Range("K19").Value = 0
$ws_CUL.Range("L19").Value = 7400.000100000001
$ws_CUL.Range("M19").ClearContents()
$ws_CUL.Range("N19").Value = -7748.000100000001

# CUL row 25
$ws_CUL.Range("H25").Value = 1747.5
$ws_CUL.Range("I25").Value = 1000
$ws_CUL.Range("J25").Value = 1996.6666
$ws_CUL.Range("K25").Value = 3000
$ws_CUL.Range("L25").Value = 5989.9998
$ws_CUL.Range("M25").Value = -2831
$ws_CUL.Range("N25").Value = -6327.9998

# CUL row 30
$ws_CUL.Range("H30").Value = 1747.5
$ws_CUL.Range("I30").Value = 1000
$ws_CUL.Range("J30").Value = 1996.6666
$ws_CUL.Range("K30").Value = 3000
$ws_CUL.Range("L30").Value = 5989.9998
$ws_CUL.Range("M30").Value = -2898
$ws_CUL.Range("N30").Value = -6193.9998

# CUL row 122
$ws_CUL.Range("H122").Value = 699.64703
$ws_CUL.Range("J122").Value = 953.7778
$ws_CUL.Range("L122").Value = 8584.0002
$ws_CUL.Range("N122").Value = -13484.0002

# CUL row 131
$ws_CUL.Range("H131").Value = 2105.1355
$ws_CUL.Range("I131").Value = 677.7778
$ws_CUL.Range("J131").Value = 2252.7932
$ws_CUL.Range("K131").Value = 2033.3334
$ws_CUL.Range("L131").Value = 6758.3796
$ws_CUL.Range("M131").Value = 3006.6666
$ws_CUL.Range("N131").Value = -16838.3796

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 59
$ws_GSM.Range("H59").Value = 25000
$ws_GSM.Range("J59").Value = 25000
$ws_GSM.Range("L59").Value = 25000
$ws_GSM.Range("N59").Value = -26166

# GSM row 122
$ws_GSM.Range("H122").Value = 2206
$ws_GSM.Range("I122").Value = 1907.2069
$ws_GSM.Range("J122").Value = 2993.7273
$ws_GSM.Range("K122").Value = 5721.620699999999
$ws_GSM.Range("L122").Value = 8981.1819
$ws_GSM.Range("M122").Value = -3271.620699999999
$ws_GSM.Range("N122").Value = -13881.1819

# GSM row 126
$ws_GSM.Range("H126").Value = 1968.7428
$ws_GSM.Range("I126").Value = 1600.625
$ws_GSM.Range("J126").Value = 2771.9092
$ws_GSM.Range("K126").Value = 4801.875
$ws_GSM.Range("L126").Value = 8315.7276
$ws_GSM.Range("M126").Value = -2331.875
$ws_GSM.Range("N126").Value = -13255.7276

# GSM row 132
$ws_GSM.Range("H132").Value = 4163.905
$ws_GSM.Range("I132").Value = 2995.625
$ws_GSM.Range("J132").Value = 7902.4
$ws_GSM.Range("K132").Value = 8986.875
$ws_GSM.Range("L132").Value = 23707.2
$ws_GSM.Range("M132").Value = -6456.875
$ws_GSM.Range("N132").Value = -28767.2

# GSM row 138
$ws_GSM.Range("H138").Value = 68730
$ws_GSM.Range("J138").Value = 68730
$ws_GSM.Range("L138").Value = 68730
$ws_GSM.Range("N138").Value = -79010

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 145186.28
$ws_LTW.Range("I7").Value = 250750
$ws_LTW.Range("J7").Value = 4434.6665
$ws_LTW.Range("K7").Value = 250750
$ws_LTW.Range("L7").Value = 4434.6665
$ws_LTW.Range("M7").Value = -250638
$ws_LTW.Range("N7").Value = -4658.6665

# LTW row 40
$ws_LTW.Range("H40").Value = 40662.152
$ws_LTW.Range("I40").Value = 51260.8
$ws_LTW.Range("J40").Value = 5333.3335
$ws_LTW.Range("K40").Value = 51260.8
$ws_LTW.Range("L40").Value = 5333.3335
$ws_LTW.Range("M40").Value = -51124.8
$ws_LTW.Range("N40").Value = -5605.3335

# LTW row 55
$ws_LTW.Range("H55").Value = 425.89474
$ws_LTW.Range("I55").Value = 442.85715
$ws_LTW.Range("J55").Value = 378.4
$ws_LTW.Range("K55").Value = 442.85715
$ws_LTW.Range("L55").Value = 378.4
$ws_LTW.Range("M55").Value = -269.85715
$ws_LTW.Range("N55").Value = -724.4

# LTW row 122
$ws_LTW.Range("H122").Value = 6540011
$ws_LTW.Range("I122").Value = 18520032
$ws_LTW.Range("K122").Value = 55560096
$ws_LTW.Range("M122").Value = -55557646

# LTW row 126
$ws_LTW.Range("H126").Value = 145186.28
$ws_LTW.Range("I126").Value = 250750
$ws_LTW.Range("J126").Value = 4434.6665
$ws_LTW.Range("K126").Value = 752250
$ws_LTW.Range("L126").Value = 13303.9995
$ws_LTW.Range("M126").Value = -749780
$ws_LTW.Range("N126").Value = -18243.9995

# LTW row 136
$ws_LTW.Range("H136").Value = 14497062
$ws_LTW.Range("I136").Value = 2806.4
$ws_LTW.Range("J136").Value = 25646488
$ws_LTW.Range("K136").Value = 8419.200000000001
$ws_LTW.Range("L136").Value = 76939464
$ws_LTW.Range("M136").Value = -5869.200000000001
$ws_LTW.Range("N136").Value = -76944564

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws_WVR.Range("H81").Value = 904.55554
$ws_WVR.Range("I81").Value = 785.25
$ws_WVR.Range("J81").Value = 1000
$ws_WVR.Range("K81").Value = 1570.5
$ws_WVR.Range("L81").Value = 2000
$ws_WVR.Range("M81").Value = -509.5
$ws_WVR.Range("N81").Value = -4122

# WVR row 84
$ws_WVR.Range("H84").Value = 904.55554
$ws_WVR.Range("I84").Value = 785.25
$ws_WVR.Range("J84").Value = 1000
$ws_WVR.Range("K84").Value = 7852.5
$ws_WVR.Range("L84").Value = 10000
$ws_WVR.Range("M84").Value = -2548.5
$ws_WVR.Range("N84").Value = -20608

# WVR row 132
$ws_WVR.Range("H132").Value = 2941.2917
$ws_WVR.Range("I132").Value = 2609.3635
$ws_WVR.Range("J132").Value = 3222.1538
$ws_WVR.Range("K132").Value = 7828.0905
$ws_WVR.Range("L132").Value = 9666.4614
$ws_WVR.Range("M132").Value = -5298.0905
$ws_WVR.Range("N132").Value = -14726.4614
